$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.926.47'
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").Value = '1.634.43'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.68'
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.505'
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.255'
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0633'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.67'
$ws.Range("E10").Value = '  +0.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0791'
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("D12").Value = '1.860.16'
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.23'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").Value = '1.650.62'
$ws.Range("E14").Value = '  +1.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.545'
$ws.Range("E15").Value = '  -1.67%  '
$ws.Range("D16").Value = '0.0₃0756'
$ws.Range("E16").Value = '  -0.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.86'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = '25.926.77'
$ws.Range("E18").Value = '  +0.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.00'
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.16'
$ws.Range("E20").Value = '  +0.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.38'
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.97'
$ws.Range("E22").Value = '  +0.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.26'
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.80'
$ws.Range("E24").Value = '  -1.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.87'
$ws.Range("E26").Value = '  +0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.126'
$ws.Range("E27").Value = '  +1.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.86'
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.48'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.24'
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0498'
$ws.Range("E31").Value = '  +0.72%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.31'
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.23'
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.58'
$ws.Range("E34").Value = '  -0.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.43'
$ws.Range("E35").Value = '  +2.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.902'
$ws.Range("E36").Value = '  -0.55%  '
$ws.Range("D37").Value = '1.140.60'
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.550'
$ws.Range("E38").Value = '  +1.38%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.48'
$ws.Range("E39").Value = '  -1.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0156'
$ws.Range("E40").Value = '  +0.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.804'
$ws.Range("E42").Value = '  +0.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.46'
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.16'
$ws.Range("E44").Value = '  -1.62%  '
$ws.Range("D45").Value = '1.769.55'
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("D46").Value = '0.0₆0112'
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.29'
$ws.Range("E47").Value = '  +2.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0526'
$ws.Range("E48").Value = '  +2.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.46'
$ws.Range("E49").Value = '  +1.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.415'
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.63'
$ws.Range("E51").Value = '  +2.16%  '
